# Fix Technology Stack section on the Enterprise Landing Zone "Solution
# Overview" slide (slide 4) to max 3 sub-items:
#   - "Cloud Platform: Microsoft Azure" ->
#     "Platform: Microsoft Azure with Management Groups and Policy framework"
#   - "Governance: Azure Policy, Azure Blueprints, Azure Resource Graph" ->
#     "Governance: Azure Policy, Blueprints, Resource Graph, RBAC"
#   - "Identity: Azure Active Directory with RBAC" paragraph removed entirely

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

function Replace-Snippet($textRange, $oldText, $newText) {
    $full = $textRange.Text
    $idx = $full.IndexOf($oldText)
    if ($idx -lt 0) {
        throw "Could not find text: $oldText"
    }
    $sub = $textRange.Characters($idx + 1, $oldText.Length)
    $sub.Text = $newText
}

function Remove-Paragraph($textRange, $paraText) {
    $full = $textRange.Text
    $idx = $full.IndexOf($paraText)
    if ($idx -lt 0) {
        throw "Could not find paragraph text: $paraText"
    }
    $len = $paraText.Length
    # Also swallow the trailing paragraph-mark (\r) so the whole paragraph
    # node is removed rather than leaving an empty paragraph behind.
    if (($idx + $len) -lt $full.Length -and $full.Substring($idx + $len, 1) -eq "`r") {
        $len = $len + 1
    }
    $sub = $textRange.Characters($idx + 1, $len)
    $sub.Delete()
}

Replace-Snippet $tr "Cloud Platform: Microsoft Azure" "Platform: Microsoft Azure with Management Groups and Policy framework"
Replace-Snippet $tr "Governance: Azure Policy, Azure Blueprints, Azure Resource Graph" "Governance: Azure Policy, Blueprints, Resource Graph, RBAC"
Remove-Paragraph $tr "Identity: Azure Active Directory with RBAC"
